# Updated WTI time series
# Applies new calibration values to start_price, Linear, and NonLinear sheets.

$wb = $excel.ActiveWorkbook

# --- start_price sheet ---
$wsStart = $wb.Worksheets.Item("start_price")
$wsStart.Range("A2").Value = 33.25

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.001211039053952787
$wsLinear.Range("B3").Value = -0.0716874430586938
$wsLinear.Range("B4").Value = 1.540928548617921
$wsLinear.Range("B5").Value = "[1.0, 0.23318776677716094, 0.03051666646448454, -0.022756364554546334, -0.02645397937035928, 0.04540380859552315, 0.2101209650488192, 0.3820628747763741, 0.20622912574808683, 0.03193087621973779, -0.018657461627024616, -0.038987960532592035, 0.036515309396303174, 0.20152390433027967, 0.3614415503858696, 0.18935763324900232, 0.020665030860140585, -0.023942182871921198, -0.024133943426347126, 0.018422551438419407]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B3").Value = 0.8988146863255276
$wsNonLinear.Range("B4").Value = -0.0007598249328794014
$wsNonLinear.Range("B5").Value = -0.03678994382682052
$wsNonLinear.Range("B6").Value = 1.709036458908619
$wsNonLinear.Range("B7").Value = 0.07264962439209825
$wsNonLinear.Range("B8").Value = -0.2264823969784011
$wsNonLinear.Range("B9").Value = 1.387586451417342
$wsNonLinear.Range("B10").Value = "[1.0, 0.23193874869759637, 0.031058775782711996, -0.02113182890330981, -0.025724106026786743, 0.04641307292125169, 0.2083310845857979, 0.38185545492132844, 0.20397538160480874, 0.03186191153099442, -0.017795622021419544, -0.03812813528153138, 0.03617071025505917, 0.20066402204449194, 0.3605414041616999, 0.188599112711134, 0.021087930260265775, -0.0233497640932082, -0.021832141684622147, 0.01806634837486836]"
